# Update the "Estadisticos" for the morning group (Matutino) statistics
# on sheets "1er Parcial" and "3er Parcial". Rows 2-5 and 10 (groups 1DV,
# 1BV, 1CV, 1EV, 1AV) get their grading results filled in, mirroring the
# already-populated rows 6-9/11-12.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

# Values: row => E(Aprobados), F(Reprobados), G(Por_Apro), H(Por_Repro), I(Promedio), J(Blancos), K(Por_Blan)
$updates = @{
    2  = @(18, 10, 64.29000000000001, 35.71, 7.8, 1, 3.57)
    3  = @(18, 26, 40.91, 59.09, 8.800000000000001, 26, 59.09)
    4  = @(17, 26, 39.53, 60.47, 8.300000000000001, 26, 60.47)
    5  = @(9, 9, 50, 50, 8.800000000000001, 9, 50)
    10 = @(24, 0, 100, 0, 8.800000000000001, 0, 0)
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]

        $ws.Cells.Item($row, 5).Value  = $vals[0]  # E - Aprobados
        $ws.Cells.Item($row, 6).Value  = $vals[1]  # F - Reprobados
        $ws.Cells.Item($row, 7).Value  = $vals[2]  # G - Por_Apro
        $ws.Cells.Item($row, 8).Value  = $vals[3]  # H - Por_Repro
        $ws.Cells.Item($row, 9).Value  = $vals[4]  # I - Promedio
        $ws.Cells.Item($row, 10).Value = $vals[5]  # J - Blancos
        $ws.Cells.Item($row, 11).Value = $vals[6]  # K - Por_Blan
    }
}
